# Auto-generated Excel COM-interop script applying numeric updates
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Ultima_Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1992.8572
$ws.Range("I112").Value = 1066.6666
$ws.Range("J112").Value = 2147.2222
$ws.Range("K112").Value = 3199.9998
$ws.Range("L112").Value = 6441.6666
$ws.Range("M112").Value = -2091.9998
$ws.Range("N112").Value = -8657.6666
# Row 125
$ws.Range("H125").Value = 2095.5715
$ws.Range("I125").Value = 2586.2
$ws.Range("J125").Value = 869
$ws.Range("K125").Value = 23275.8
$ws.Range("L125").Value = 7821
$ws.Range("M125").Value = -20815.8
$ws.Range("N125").Value = -12741
# Row 138
$ws.Range("H138").Value = 3307.6565
$ws.Range("I138").Value = 1574.4736
$ws.Range("J138").Value = 4387.344
$ws.Range("K138").Value = 4723.4208
$ws.Range("L138").Value = 13162.032
$ws.Range("M138").Value = 416.5792000000001
$ws.Range("N138").Value = -23442.032

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1812.2858
$ws.Range("I61").Value = 1723.4468
$ws.Range("J61").Value = 3900
$ws.Range("K61").Value = 1723.4468
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -1511.4468
$ws.Range("N61").Value = -4324
# Row 136
$ws.Range("H136").Value = 1812.2858
$ws.Range("I136").Value = 1723.4468
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 5170.3404
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -2620.3404
$ws.Range("N136").Value = -16800

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1085.2727
$ws.Range("I99").Value = 577.1429000000001
$ws.Range("J99").Value = 1974.5
$ws.Range("K99").Value = 577.1429000000001
$ws.Range("L99").Value = 1974.5
$ws.Range("M99").Value = 920.8570999999999
$ws.Range("N99").Value = -4970.5
# Row 118
$ws.Range("H118").Value = 7877.273
$ws.Range("J118").Value = 7877.273
$ws.Range("L118").Value = 7877.273
$ws.Range("N118").Value = -11191.273
# Row 134
$ws.Range("H134").Value = 4048.7104
$ws.Range("I134").Value = 3090.6128
$ws.Range("J134").Value = 8291.714
$ws.Range("K134").Value = 9271.838400000001
$ws.Range("L134").Value = 24875.142
$ws.Range("M134").Value = -6736.838400000001
$ws.Range("N134").Value = -29945.142

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 7812.467
$ws.Range("I16").Value = 11304.556
$ws.Range("J16").Value = 2574.3333
$ws.Range("K16").Value = 11304.556
$ws.Range("L16").Value = 2574.3333
$ws.Range("M16").Value = -11017.556
$ws.Range("N16").Value = -3148.3333
# Row 31
$ws.Range("H31").Value = 5212.4062
$ws.Range("I31").Value = 1144.08
$ws.Range("J31").Value = 7820.3076
$ws.Range("K31").Value = 1144.08
$ws.Range("L31").Value = 7820.3076
$ws.Range("M31").Value = -849.0799999999999
$ws.Range("N31").Value = -8410.3076
# Row 34
$ws.Range("H34").Value = 5212.4062
$ws.Range("I34").Value = 1144.08
$ws.Range("J34").Value = 7820.3076
$ws.Range("K34").Value = 1144.08
$ws.Range("L34").Value = 7820.3076
$ws.Range("M34").Value = -942.0799999999999
$ws.Range("N34").Value = -8224.3076
# Row 99
$ws.Range("H99").Value = 1700
$ws.Range("I99").Value = 1700
$ws.Range("K99").Value = 1700
$ws.Range("M99").Value = -202
# Row 113
$ws.Range("H113").Value = 7812.467
$ws.Range("I113").Value = 11304.556
$ws.Range("J113").Value = 2574.3333
$ws.Range("K113").Value = 11304.556
$ws.Range("L113").Value = 2574.3333
$ws.Range("M113").Value = -9134.556
$ws.Range("N113").Value = -6914.3333
# Row 126
$ws.Range("H126").Value = 1700
$ws.Range("I126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("M126").Value = -2630
# Row 134
$ws.Range("H134").Value = 883503.5600000001
$ws.Range("I134").Value = 1928.6
$ws.Range("J134").Value = 3402289.2
$ws.Range("K134").Value = 5785.799999999999
$ws.Range("L134").Value = 10206867.6
$ws.Range("M134").Value = -3250.799999999999
$ws.Range("N134").Value = -10211937.6

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1866.8235
$ws.Range("I132").Value = 1302.3334
$ws.Range("J132").Value = 1987.7858
$ws.Range("K132").Value = 11721.0006
$ws.Range("L132").Value = 17890.0722
$ws.Range("M132").Value = -9191.000599999999
$ws.Range("N132").Value = -22950.0722

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 4899.5454
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4899.5454
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 14698.6362
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -19638.6362

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7139.3
$ws.Range("I7").Value = 5377.2144
$ws.Range("J7").Value = 8681.125
$ws.Range("K7").Value = 5377.2144
$ws.Range("L7").Value = 8681.125
$ws.Range("M7").Value = -5265.2144
$ws.Range("N7").Value = -8905.125
# Row 55
$ws.Range("H55").Value = 353.2
$ws.Range("I55").Value = 207.1
$ws.Range("J55").Value = 499.3
$ws.Range("K55").Value = 207.1
$ws.Range("L55").Value = 499.3
$ws.Range("M55").Value = -34.09999999999999
$ws.Range("N55").Value = -845.3
# Row 63
$ws.Range("H63").Value = 18018.334
$ws.Range("J63").Value = 18018.334
$ws.Range("L63").Value = 18018.334
$ws.Range("N63").Value = -19516.334
# Row 66
$ws.Range("H66").Value = 18018.334
$ws.Range("J66").Value = 18018.334
$ws.Range("L66").Value = 54055.00199999999
$ws.Range("N66").Value = -61543.00199999999
# Row 68
$ws.Range("H68").Value = 2502.6956
$ws.Range("I68").Value = 2042.2778
$ws.Range("K68").Value = 2042.2778
$ws.Range("M68").Value = -1293.2778
# Row 71
$ws.Range("H71").Value = 2502.6956
$ws.Range("I71").Value = 2042.2778
$ws.Range("K71").Value = 10211.389
$ws.Range("M71").Value = -6467.389000000001
# Row 93
$ws.Range("H93").Value = 1327.3704
$ws.Range("I93").Value = 1291.8422
$ws.Range("J93").Value = 1411.75
$ws.Range("K93").Value = 1291.8422
$ws.Range("L93").Value = 1411.75
$ws.Range("M93").Value = -43.84220000000005
$ws.Range("N93").Value = -3907.75
# Row 126
$ws.Range("H126").Value = 7139.3
$ws.Range("I126").Value = 5377.2144
$ws.Range("J126").Value = 8681.125
$ws.Range("K126").Value = 16131.6432
$ws.Range("L126").Value = 26043.375
$ws.Range("M126").Value = -13661.6432
$ws.Range("N126").Value = -30983.375

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 641
$ws.Range("I113").Value = 702
$ws.Range("J113").Value = 580
$ws.Range("K113").Value = 2106
$ws.Range("L113").Value = 1740
$ws.Range("M113").Value = 64
$ws.Range("N113").Value = -6080
# Row 126
$ws.Range("H126").Value = 2000.439
$ws.Range("I126").Value = 2037.5264
$ws.Range("K126").Value = 6112.5792
$ws.Range("M126").Value = -3642.5792
